$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.054.01"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.564.26"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0598"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.787.78"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "1.563.16"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "27.051.57"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "0.0₃0706"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("D34").Value = "1.424.20"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.35%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.533"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "1.701.39"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0518"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
